$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.876.13"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").Value = "'3.256.94"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'583.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("D6").Value = "'183.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.73%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("E9").Value = '  +4.51%  '
$ws.Range("D10").Value = "'6.67"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = "'3.822.48"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = "'28.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.68%  '
$ws.Range("D15").Value = "'67.885.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("E16").Value = '  +2.84%  '
$ws.Range("D17").Value = "'3.257.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = "'5.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = "'13.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("D20").Value = "'381.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.41%  '
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").Value = "'71.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = "'0.0000120"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("D26").Value = "'9.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").Value = "'22.91"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("D32").Value = "'7.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.51%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  +2.81%  '
$ws.Range("E35").Value = '  +2.75%  '
$ws.Range("D36").Value = "'161.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.04%  '
$ws.Range("D38").Value = "'0.835"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("D39").Value = "'26.63"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("E40").Value = '  +7.61%  '
$ws.Range("D41").Value = "'6.67"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").Value = "'2.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = "'41.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").Value = "'347.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("D45").Value = "'25.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.43%  '
$ws.Range("D46").Value = "'0.0688"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("D47").Value = "'2.634.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("E51").Value = '  +3.25%  '
